$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "01-07-2021"
$ws.Range("A16").Style = "Normal"
$ws.Range("B16").Value = 110.45
$ws.Range("C16").Value = 108.62
$ws.Range("D16").Value = 112.15
$ws.Range("E16").Value = 108.45
$ws.Range("F16").Value = 119.48
